$d = $word.ActiveDocument

$d.Content.Find.Execute("kjkim761@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "kjkim761@proton.me", 2)

$d.PageSetup.Orientation = 0
